# Applies the two content edits described by the commit:
#  1. Move the "_GoBack" bookmark from the very start of the document
#     (the title paragraph) down to the last (empty) paragraph of the
#     body, just before the sectPr.
#  2. Re-run a "spelling/grammar check" pass over the
#     "StreetWise Partners, Inc., Washington, D.C. " job-location line,
#     which splits the two bold runs into several runs and brackets
#     "StreetWise" / "Washington, D.C." with <w:proofErr/> markers
#     (spellStart/spellEnd around "StreetWise", gramStart/gramEnd around
#     "Washington, D.C."). The visible text is unchanged.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Move the _GoBack bookmark to the end of the document.
# ---------------------------------------------------------------------
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

$lastPara = $d.Paragraphs.Last.Range
$d.Bookmarks.Add("_GoBack", $lastPara)

# ---------------------------------------------------------------------
# 2. Split the StreetWise/Washington runs and add proofErr markers.
# ---------------------------------------------------------------------
$target = $d.Content
$found = $target.Find.Execute("StreetWise Partners, Inc., Washington, D.C. ", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertAt = $target.Start
    $target.Delete()

    $xml = @'
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:proofErr w:type="spellStart"/>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
              <w:t>StreetWise</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
              <w:t xml:space="preserve"> Partners, Inc., </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
              <w:t>Washington, D.C.</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:sz w:val="20"/>
                <w:szCs w:val="20"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

    $insertRange = $d.Range($insertAt, $insertAt)
    $insertRange.InsertXML($xml)
}

Write-Output "done"
